$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65-184 down to 66-185.
$ws.Rows(65).Insert()

# Populate the newly inserted row 65 with the new data record.
$ws.Cells.Item(65, 1).Value = 4
$ws.Cells.Item(65, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(65, 3).Value = "Los Lagos"
$ws.Cells.Item(65, 4).Value = 44536
$ws.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(65, 5).Value = 10
$ws.Cells.Item(65, 6).Value = "Fruta"
$ws.Cells.Item(65, 7).Value = 100102
$ws.Cells.Item(65, 8).Value = "Cítricos"
$ws.Cells.Item(65, 9).Value = 100102006
$ws.Cells.Item(65, 10).Value = "Pomelo"
$ws.Cells.Item(65, 11).Value = "Start Ruby"
$ws.Cells.Item(65, 12).Value = "Primera"
$ws.Cells.Item(65, 13).Value = 120
$ws.Cells.Item(65, 14).Value = 11000
$ws.Cells.Item(65, 15).Value = 12000
$ws.Cells.Item(65, 16).Value = 11500
$ws.Cells.Item(65, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(65, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(65, 19).Value = 821
$ws.Cells.Item(65, 20).Value = 14
